$d = $word.ActiveDocument

# 1. Replace text about requesting appointments (remove "y rector")
$d.Content.Find.Execute(
    "Puede solicitar citas de atención a profesores, directores académicos y rector.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Puede solicitar citas de atención a profesores y directores académicos.", 2)

# 2. The paragraph that used to talk about "rechazar" now talks about "aceptar"
$d.Content.Find.Execute(
    "Puede rechazar una fecha y hora sugerida por un profesor, director académico o rector.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Puede aceptar la fecha y hora sugeridas por un profesor, director académico.", 2)

# 3. The paragraph that used to talk about "aceptar" now talks about "rechazar"
$d.Content.Find.Execute(
    "Puede aceptar la fecha y hora sugeridas por un profesor, director académico o rector.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Puede rechazar una fecha y hora sugerida por un profesor, director académico.", 2)

# 4. Replace text about evaluating appointments (remove "directores académicos y rector")
$d.Content.Find.Execute(
    "Puede evaluar citas con profesores, directores académicos y rector.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Puede evaluar citas con profesores.", 2)

# 5. Remove the two paragraphs "Puede crear carreras." and "Puede modificar carreras."
foreach ($text in @("Puede crear carreras.", "Puede modificar carreras.")) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.TrimEnd("`r`a") -eq $text) {
            $p.Range.Delete()
            break
        }
    }
}
